$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "rajasri1213"
$ws.Range("B2").Value = "rajasr13423i@gmail.com"
$ws.Range("A3").Value = "vineela12563"
$ws.Range("B3").Value = "vineela125673@gmail.com"

$ws.Range("B3").Select()
